$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reviewer comment text that used to live in E4/F4 was moved out of the
# spreadsheet (converted into a separate Word document per the commit
# message), so clear those two cells.
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# Restore their plain (unfilled) Arial 14 formatting explicitly so the
# underlying style record matches (font reused, no wrap, no fill).
$r = $ws.Range("E4:F4")
$r.Font.Name = "Arial"
$r.Font.Size = 14
$r.Font.Color = 2367776

# A11/C11/D11/D10 likewise lose their "fill applied" styling remnants; copy
# the plain format already used by equivalent cells elsewhere in the sheet.
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C4").Copy()
$ws.Range("C11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

# Update the view/selection to match the edited workbook.
$ws.Range("F4").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 4
